$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 171.09091
$ws.Range("I12").Value = 157.14285
$ws.Range("J12").Value = 195.5
$ws.Range("K12").Value = 157.14285
$ws.Range("L12").Value = 195.5
$ws.Range("M12").Value = 12.85714999999999
$ws.Range("N12").Value = -535.5

$ws.Range("H33").Value = 461.57144
$ws.Range("I33").Value = 390
$ws.Range("J33").Value = 590.4
$ws.Range("K33").Value = 390
$ws.Range("L33").Value = 590.4
$ws.Range("M33").Value = -161
$ws.Range("N33").Value = -1048.4

$ws.Range("H62").Value = 15637357
$ws.Range("I62").Value = 19244000
$ws.Range("J62").Value = 8566.666999999999
$ws.Range("K62").Value = 19244000
$ws.Range("L62").Value = 8566.666999999999
$ws.Range("M62").Value = -19243376
$ws.Range("N62").Value = -9814.666999999999

$ws.Range("H65").Value = 15637357
$ws.Range("I65").Value = 19244000
$ws.Range("J65").Value = 8566.666999999999
$ws.Range("K65").Value = 96220000
$ws.Range("L65").Value = 42833.335
$ws.Range("M65").Value = -96216880
$ws.Range("N65").Value = -49073.335

$ws.Range("H86").Value = 17983.334
$ws.Range("I86").Value = 34466.668
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 34466.668
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -33343.668
$ws.Range("N86").Value = -3746

$ws.Range("H89").Value = 17983.334
$ws.Range("I89").Value = 34466.668
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 172333.34
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -166717.34
$ws.Range("N89").Value = -18732

$ws.Range("H98").Value = 156251260
$ws.Range("I98").Value = 250000400
$ws.Range("J98").Value = 2693.6667
$ws.Range("K98").Value = 250000400
$ws.Range("L98").Value = 2693.6667
$ws.Range("M98").Value = -249998902
$ws.Range("N98").Value = -5689.6667

$ws.Range("H122").Value = 156251260
$ws.Range("I122").Value = 250000400
$ws.Range("J122").Value = 2693.6667
$ws.Range("K122").Value = 750001200
$ws.Range("L122").Value = 8081.000100000001
$ws.Range("M122").Value = -749998750
$ws.Range("N122").Value = -12981.0001

$ws.Range("H138").Value = 3448.4895
$ws.Range("I138").Value = 2159.5925
$ws.Range("J138").Value = 3952.8406
$ws.Range("K138").Value = 6478.7775
$ws.Range("L138").Value = 11858.5218
$ws.Range("M138").Value = -1338.7775
$ws.Range("N138").Value = -22138.5218

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12670928
$ws.Range("I32").Value = 13709292
$ws.Range("J32").Value = 37499.5
$ws.Range("K32").Value = 13709292
$ws.Range("L32").Value = 37499.5
$ws.Range("M32").Value = -13709005
$ws.Range("N32").Value = -38073.5

$ws.Range("H97").Value = 1163.8928
$ws.Range("I97").Value = 1007.7857
$ws.Range("J97").Value = 1320
$ws.Range("K97").Value = 1007.7857
$ws.Range("L97").Value = 1320
$ws.Range("M97").Value = -511.7857
$ws.Range("N97").Value = -2312

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 5448.3213
$ws.Range("I80").Value = 2333.6667
$ws.Range("J80").Value = 7784.3125
$ws.Range("K80").Value = 2333.6667
$ws.Range("L80").Value = 7784.3125
$ws.Range("M80").Value = -1335.6667
$ws.Range("N80").Value = -9780.3125

$ws.Range("H83").Value = 5448.3213
$ws.Range("I83").Value = 2333.6667
$ws.Range("J83").Value = 7784.3125
$ws.Range("K83").Value = 11668.3335
$ws.Range("L83").Value = 38921.5625
$ws.Range("M83").Value = -6676.333500000001
$ws.Range("N83").Value = -48905.5625

$ws.Range("H86").Value = 1293570.1
$ws.Range("I86").Value = 2051.5715
$ws.Range("J86").Value = 2115445.8
$ws.Range("K86").Value = 2051.5715
$ws.Range("L86").Value = 2115445.8
$ws.Range("M86").Value = -928.5715
$ws.Range("N86").Value = -2117691.8

$ws.Range("H89").Value = 1293570.1
$ws.Range("I89").Value = 2051.5715
$ws.Range("J89").Value = 2115445.8
$ws.Range("K89").Value = 10257.8575
$ws.Range("L89").Value = 10577229
$ws.Range("M89").Value = -4641.8575
$ws.Range("N89").Value = -10588461

$ws.Range("H107").Value = 33334014
$ws.Range("I107").Value = 35714908
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 35714908
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -35712988
$ws.Range("N107").Value = -5340

$ws.Range("H134").Value = 3713418.8
$ws.Range("I134").Value = 8532.471
$ws.Range("J134").Value = 8558270
$ws.Range("K134").Value = 25597.413
$ws.Range("L134").Value = 25674810
$ws.Range("M134").Value = -23062.413
$ws.Range("N134").Value = -25679880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3186.2444
$ws.Range("I31").Value = 4021.5
$ws.Range("J31").Value = 2947.6
$ws.Range("K31").Value = 4021.5
$ws.Range("L31").Value = 2947.6
$ws.Range("M31").Value = -3726.5
$ws.Range("N31").Value = -3537.6

$ws.Range("H34").Value = 3186.2444
$ws.Range("I34").Value = 4021.5
$ws.Range("J34").Value = 2947.6
$ws.Range("K34").Value = 4021.5
$ws.Range("L34").Value = 2947.6
$ws.Range("M34").Value = -3819.5
$ws.Range("N34").Value = -3351.6

$ws.Range("H107").Value = 576.7619
$ws.Range("I107").Value = 394.5
$ws.Range("J107").Value = 1160
$ws.Range("K107").Value = 394.5
$ws.Range("L107").Value = 1160
$ws.Range("M107").Value = 1525.5
$ws.Range("N107").Value = -5000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 23812208
$ws.Range("I5").Value = 5747541.5
$ws.Range("J5").Value = 50005976
$ws.Range("K5").Value = 17242624.5
$ws.Range("L5").Value = 150017928
$ws.Range("M5").Value = -17242512.5
$ws.Range("N5").Value = -150018152

$ws.Range("H49").Value = 1125
$ws.Range("J49").Value = 1125
$ws.Range("L49").Value = 3375
$ws.Range("N49").Value = -3687

$ws.Range("H131").Value = 891.21
$ws.Range("J131").Value = 930.337
$ws.Range("L131").Value = 2791.011
$ws.Range("N131").Value = -12871.011

$ws.Range("H135").Value = 23812208
$ws.Range("I135").Value = 5747541.5
$ws.Range("J135").Value = 50005976
$ws.Range("K135").Value = 51727873.5
$ws.Range("L135").Value = 450053784
$ws.Range("M135").Value = -51725338.5
$ws.Range("N135").Value = -450058854

$ws.Range("H137").Value = 2748.4856
$ws.Range("I137").Value = 2223.125
$ws.Range("J137").Value = 3190.8948
$ws.Range("K137").Value = 6669.375
$ws.Range("L137").Value = 9572.6844
$ws.Range("M137").Value = -1569.375
$ws.Range("N137").Value = -19772.6844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 986.4211
$ws.Range("I102").Value = 824.9091
$ws.Range("J102").Value = 1208.5
$ws.Range("K102").Value = 824.9091
$ws.Range("L102").Value = 1208.5
$ws.Range("M102").Value = 797.0909
$ws.Range("N102").Value = -4452.5

$ws.Range("H126").Value = 3796.25
$ws.Range("I126").Value = 3110
$ws.Range("J126").Value = 4482.5
$ws.Range("K126").Value = 9330
$ws.Range("L126").Value = 13447.5
$ws.Range("M126").Value = -6860
$ws.Range("N126").Value = -18387.5

$ws.Range("H132").Value = 5602.207
$ws.Range("I132").Value = 1493.2632
$ws.Range("J132").Value = 13409.2
$ws.Range("K132").Value = 4479.7896
$ws.Range("L132").Value = 40227.60000000001
$ws.Range("M132").Value = -1949.7896
$ws.Range("N132").Value = -45287.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5953037.5
$ws.Range("I46").Value = 8333952.5
$ws.Range("J46").Value = 750
$ws.Range("K46").Value = 8333952.5
$ws.Range("L46").Value = 750
$ws.Range("M46").Value = -8333764.5
$ws.Range("N46").Value = -1126

$ws.Range("H55").Value = 125012660
$ws.Range("I55").Value = 50115.5
$ws.Range("J55").Value = 166666830
$ws.Range("K55").Value = 50115.5
$ws.Range("L55").Value = 166666830
$ws.Range("M55").Value = -49942.5
$ws.Range("N55").Value = -166667176

$ws.Range("H132").Value = 15630015
$ws.Range("I132").Value = 33335206
$ws.Range("J132").Value = 7788.4707
$ws.Range("K132").Value = 100005618
$ws.Range("L132").Value = 23365.4121
$ws.Range("M132").Value = -100003088
$ws.Range("N132").Value = -28425.4121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 41667230
$ws.Range("I81").Value = 55556150
$ws.Range("J81").Value = 463
$ws.Range("K81").Value = 111112300
$ws.Range("L81").Value = 926
$ws.Range("M81").Value = -111111239
$ws.Range("N81").Value = -3048

$ws.Range("H84").Value = 41667230
$ws.Range("I84").Value = 55556150
$ws.Range("J84").Value = 463
$ws.Range("K84").Value = 555561500
$ws.Range("L84").Value = 4630
$ws.Range("M84").Value = -555556196
$ws.Range("N84").Value = -15238

$ws.Range("H136").Value = 3280.7964
$ws.Range("I136").Value = 4904
$ws.Range("J136").Value = 1657.5927
$ws.Range("K136").Value = 14712
$ws.Range("L136").Value = 4972.7781
$ws.Range("M136").Value = -12162
$ws.Range("N136").Value = -10072.7781
